$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.461700000000002
$ws.Range("C9").Value = -10.1406
$ws.Range("D12").Value = -7.079499999999997
$ws.Range("C18").Value = -12.54239999999999
$ws.Range("C20").Value = -11.7382
$ws.Range("D26").Value = -8.555800000000001
$ws.Range("C27").Value = -12.95879999999999
$ws.Range("D27").Value = -8.980699999999999
$ws.Range("D29").Value = -7.471800000000002
$ws.Range("D37").Value = -7.5768
$ws.Range("D38").Value = -8.4499
$ws.Range("D51").Value = -7.757
$ws.Range("D55").Value = -8.500799999999996
$ws.Range("C69").Value = -11.41140000000001
$ws.Range("D69").Value = -7.182999999999997
$ws.Range("D70").Value = -7.438999999999997
$ws.Range("C76").Value = -12.09070000000001
$ws.Range("C82").Value = -11.69519999999999
$ws.Range("D83").Value = -8.8255
$ws.Range("D102").Value = -7.6496
